# Rename the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Lotarea - Landuse"

# New column E header + values ("Count of Records Changed")
$ws.Range("E1").Value = "Count of Records Changed"
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 66
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 21
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 1966
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 6
$ws.Range("E11").Value = 4
$ws.Range("E12").Value = 7
$ws.Range("E13").Value = 20

# Match the numeric/thousands format used by column D for the new column
$ws.Range("E2:E13").NumberFormat = $ws.Range("D2:D13").NumberFormat

# Header row: bold + wrap text (matches existing header style, now wrapped)
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").WrapText = $true

# Column E width (~16.8 chars, matches the column width set alongside the new field)
$ws.Columns.Item(5).ColumnWidth = 16

# Row 1 taller to fit the wrapped header text
$ws.Rows.Item(1).RowHeight = 28.5

# Selection ends on C4 (matches the saved view state in the edited file)
$ws.Range("C4").Select() | Out-Null
